$wb = $excel.ActiveWorkbook

# Rename the existing data sheet and set its selection, then add the new Instructions sheet in front of it
$rain = $wb.Worksheets.Item(1)
$rain.Name = "rain_data"
$rain.Range("E27").Select()

$instructions = $wb.Worksheets.Add()
$instructions.Name = "Instructions"
$instructions.Move($rain)

# Usage guide text, one line per row in column A
$lines = @(
  'Application Usage Guide',
  '=======================',
  'This guide provides instructions for using the application and outlines the data requirements for successful analysis.',
  'Analysis Types',
  '--------------',
  'IMPORTANT: The column names and tab names in the templates (both rainfall and flow) must remain unchanged for accurate analysis and app''s functionality.',
  'The application supports two types of analysis: Rainfall Analysis and Flow Analysis. In the zip package that ',
  'Rainfall Analysis',
  '-----------------',
  'To perform a Rainfall Analysis, follow these steps:',
  '1. Open the rainfall template (rainfall_template.xlsx).',
  '2. Copy and paste the rainfall data into the template.',
  '3. Ensure that the rain gauge values are entered into the designated rain column within the template.',
  'Demo data is available for flow analysis under demo_data folder.',
  'Flow Analysis',
  '-------------',
  'To perform a Flow Analysis, follow these steps:',
  '1. Open the flow template (flow_template.xlsx).',
  '2. Copy and paste the flow data from your datasheet into the downloaded template.',
  '3. Accommodate up to four flows for analysis: inflow 1, inflow 2, outflow, and bypass.',
  '4. Refer to the Methods tab in the template for an illustration of the possible flow type configurations.',
  '5. You are not required to submit data for all four flow types; any combination of flow types is acceptable.',
  'IMPORTANT: If any of the flow data types are not applicable, leave the corresponding tab blank. Do not modify the columns or the tab names.',
  'Demo data is available for flow analysis under demo_data folder.',
  'Data Requirements',
  '-----------------',
  'To ensure successful analysis, the uploaded Excel spreadsheet must conform to the following requirements:',
  '- Flow must be reported in units of L/s, gpm, or cfs.',
  '- The timestamp should be in 24-hour format (mm/dd/yy hh:mm:ss).',
  '- Each tab must contain exactly two columns: one for the sample timestamps data and one for the associated values.',
  '- The column names and the tab names must not be changed from the template.',
  'Please ensure that your data meets these requirements before using the application for analysis.',
  'For further assistance or inquiries, please contact our support team at stormwater@sccwrp.org'
)

for ($i = 0; $i -lt $lines.Count; $i++) {
  $row = $i + 1
  $text = $lines[$i]
  # Guard against Excel parsing a leading "=" as a formula trigger
  if ($text.StartsWith("=")) {
    $text = "'" + $text
  }
  $instructions.Cells.Item($row, 1).Value = $text
}

# Match the source workbook formatting: plain (non-bold) black text, wide column
$dataRange = $instructions.Range("A1:A33")
$dataRange.Font.Color = 0
$instructions.Columns.Item(1).ColumnWidth = 143.1666666667

$instructions.Range("A33").Select()
